$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header/index row, extend from B1:D1 (4,5,6) to B1:J1 (0..8) ---
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4
$ws.Range("G1").Value = 5
$ws.Range("H1").Value = 6
$ws.Range("I1").Value = 7
$ws.Range("J1").Value = 8

# New cells E1:J1 need the same formatting that B1:D1 already carry
# (bold font, thin box border, centered / top-aligned).
$ws.Range("E1:J1").Font.Bold = $true
$ws.Range("E1:J1").HorizontalAlignment = -4108
$ws.Range("E1:J1").VerticalAlignment = -4160
$ws.Range("E1:J1").Borders.LineStyle = 1

# --- Row 2: carID ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "carID"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 6
$ws.Range("H2").Value = 9
$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 12

# --- Row 3: speed2 ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "speed2"
$ws.Range("C3").Value = 42.43
$ws.Range("D3").Value = 55.9
$ws.Range("E3").Value = 49.24
$ws.Range("F3").Value = 43.01
$ws.Range("G3").Value = 40.31
$ws.Range("H3").Value = 47.17
$ws.Range("I3").Value = 40.31
$ws.Range("J3").Value = 44.72

# --- Row 4: asma ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "asma"
$ws.Range("C4").Value = 21.22
$ws.Range("D4").Value = 59.72
$ws.Range("E4").Value = 40.7
$ws.Range("F4").Value = 22.89
$ws.Range("G4").Value = 15.18
$ws.Range("H4").Value = 34.77
$ws.Range("I4").Value = 15.18
$ws.Range("J4").Value = 27.78

# --- Row 5: ceza_tutar ---
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "ceza_tutar"
$ws.Range("C5").Value = 1508.5
$ws.Range("D5").Value = 6440
$ws.Range("E5").Value = 3136
$ws.Range("F5").Value = 1508.5
$ws.Range("G5").Value = 1508.5
$ws.Range("H5").Value = 3136
$ws.Range("I5").Value = 1508.5
$ws.Range("J5").Value = 1508.5

# --- Row 6: hesaplanan_asma ---
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "hesaplanan_asma"
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 50
$ws.Range("E6").Value = 30
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 30
$ws.Range("I6").Value = 10
$ws.Range("J6").Value = 10
